# Added DataProvider with excel
# invalidLoginTest: replace "Kevin" with "Bala" on row 2, and add a new
# row 4 (Kim / kim124 / Invalid credentials) below the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalidLoginTest")

# Clone formatting (style + row height) of row 3 into the new row 4 before
# writing values, so the new row matches the look of the existing data rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = 15.75

# Update existing row 2 username.
$ws.Range("A2").Value = "Bala"

# Populate the newly added row 4.
$ws.Range("A4").Value = "Kim"
$ws.Range("B4").Value = "kim124"
$ws.Range("C4").Value = "Invalid credentials"

# Match the author's final selection state.
[void]$ws.Range("A3:C4").Select()
